# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet, which both carry duplicated rows for the
# same events.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 157
$wsExhibit.Range("F3").Value = 472
$wsExhibit.Range("F4").Value = 12
$wsExhibit.Range("F9").Value = 263

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 157
$wsAll.Range("F4").Value = 472
$wsAll.Range("F5").Value = 12
$wsAll.Range("F10").Value = 263
